$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Reference style/format cells already present in the sheet
$dateFormat = $ws.Range("C3").NumberFormat

# Rows 19-23: "Conception maquette IHM" task renamed to "Conception maquette IHM photoshop"
$ws.Range("B19").Value = "Conception maquette IHM photoshop"
$ws.Range("B20").Value = "Conception maquette IHM photoshop"
$ws.Range("B21").Value = "Conception maquette IHM photoshop"
$ws.Range("B22").Value = "Conception maquette IHM photoshop"
$ws.Range("B23").Value = "Conception maquette IHM photoshop"

# New rows 25-28 with new tasks/dates/hours
$ws.Range("B25").Value = "Installation d'un raspberry de test"
$ws.Range("C25").Value = 42326
$ws.Range("C25").NumberFormat = $dateFormat
$ws.Range("D25").Value = 0.2

$ws.Range("B26").Value = "Actualisation suivi d'activité"
$ws.Range("C26").Value = 42326
$ws.Range("C26").NumberFormat = $dateFormat
$ws.Range("D26").Value = 0.25

$ws.Range("B27").Value = "Actualisation du gantt"
$ws.Range("C27").Value = 42326
$ws.Range("C27").NumberFormat = $dateFormat
$ws.Range("D27").Value = 0.5

$ws.Range("B28").Value = "Actualisation suivi des tâches"
$ws.Range("C28").Value = 42326
$ws.Range("C28").NumberFormat = $dateFormat
$ws.Range("D28").Value = 1

# Update the selected cell in the sheet view
$ws.Range("F19").Select()
